$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) values must stay text even though many look numeric;
# force Text number format before assignment, then reset the style index
# back to Normal so no residual style/format change is left on the cell.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.249.77"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +2.56%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.585.78"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +1.26%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "213.59"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +1.26%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "23.95"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +6.64%  "
$ws.Range("E9").Value = "  +0.13%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0597"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0886"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.20%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.813.55"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.54%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.586.64"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  -0.26%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "28.274.38"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.16"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "227.68"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("E19").Value = "  +0.04%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.46"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +1.18%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.06"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.33"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  +0.87%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.08"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -0.07%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.401.42"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("E36").Value = "  -8.35%  "
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").Value = "  -0.39%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.53"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +8.61%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.540"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.809"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  +1.18%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.88"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.57%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "5.60"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.979"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "64.23"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.723.03"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("E48").Value = "  +1.58%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "86.84"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  +4.23%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0520"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
